# Update countries & provincias Spain
# Refreshes the "Pais" sheet with a newer data pull: updates the
# "last updated" timestamp and refreshes case counts for several
# countries. A handful of countries overtook their neighbours in the
# total-cases ranking, so those rows were re-sorted to match; we apply
# that by writing the new per-row numbers directly (rank position is
# driven by column B, "Casos totales").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 10:22"

# --- Simple in-place refreshes (rank unchanged) ------------------------
Set-RowValues 17 @(11821, 40, 2998, 8619, 244, 18, 204)   # Austria
Set-RowValues 30 @(3834, 207, 134, 3616, 50, 5, 84)       # Polonia
Set-RowValues 68 @(752, 0, 29, 709, 65, 2, 14)            # Moldavia
Set-RowValues 72 @(629, 5, 30, 578, 4, 0, 21)             # Bosnia y Herzegovina

# --- Australia overtakes Noruega (rows 23-24) ---------------------------
$ws.Range("A23").Value = "Australia"
$ws.Range("A24").Value = "Noruega"
Set-RowValues 23 @(5687, 137, 585, 5068, 85, 4, 34)
Set-RowValues 24 @(5645, 95, 32, 5551, 98, 0, 62)

# --- Filipinas overtakes Japon (rows 35-36) ------------------------------
$ws.Range("A35").Value = "Filipinas"
$ws.Range("A36").Value = "Japon"
Set-RowValues 35 @(3246, 152, 64, 3030, 1, 8, 152)
Set-RowValues 36 @(3139, 0, 514, 2548, 64, 0, 77)

# --- Kazajistan overtakes Camerun & Tunez (rows 73-75) -------------------
$ws.Range("A73").Value = "Kazajistan"
$ws.Range("A74").Value = "Camerun"
$ws.Range("A75").Value = "Tunez"
Set-RowValues 73 @(569, 38, 36, 527, 6, 1, 6)
Set-RowValues 74 @(555, 0, 17, 529, 0, 0, 9)
Set-RowValues 75 @(553, 0, 5, 530, 26, 0, 18)

# --- Laos overtakes Groenlandia & Curazao (rows 175-177) -----------------
$ws.Range("A175").Value = "Laos"
$ws.Range("A176").Value = "Groenlandia"
$ws.Range("A177").Value = "Curazao"
Set-RowValues 175 @(11, 1, 0, 11, 0, 0, 0)
Set-RowValues 176 @(11, 0, 3, 8, 0, 0, 0)
Set-RowValues 177 @(11, 0, 5, 5, 0, 0, 1)
